$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "## ... ##" markers from the three phonetic-gloss cells.
$ws.Range("H5").Value = "<peteĩ>"
$ws.Range("H9").Value = "<mokõi>"
$ws.Range("H11").Value = "<mbohapy>"

# Move the active selection in the frozen bottom-right pane to K25.
$ws.Range("K25").Select()
